# Apply updated crypto price/volume data (and two name/link swaps) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.348.22'
$ws.Range('E2').Value = '  -1.30%  '
$ws.Range('D3').Value = '1.815.30'
$ws.Range('E3').Value = '  -3.35%  '
$ws.Range('D4').Value = '''1.004'
$ws.Range('E4').Value = '  -1.02%  '
$ws.Range('D5').Value = '''330.73'
$ws.Range('E5').Value = '  -1.61%  '
$ws.Range('E6').Value = '  -0.83%  '
$ws.Range('D7').Value = '''0.4553'
$ws.Range('E7').Value = '  -2.11%  '
$ws.Range('D8').Value = '''0.3799'
$ws.Range('E8').Value = '  -3.67%  '
$ws.Range('E9').Value = '  +0.16%  '
$ws.Range('D10').Value = '''0.07824'
$ws.Range('E10').Value = '  -2.17%  '
$ws.Range('D11').Value = '''0.9583'
$ws.Range('E11').Value = '  -4.98%  '
$ws.Range('D12').Value = '''20.90'
$ws.Range('E12').Value = '  -4.62%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '''5.832'
$ws.Range('E13').Value = '  -2.56%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.808.35'
$ws.Range('E14').Value = '  -4.06%  '
$ws.Range('D15').Value = '''7.051'
$ws.Range('E15').Value = '  -2.93%  '
$ws.Range('E16').Value = '  -1.09%  '
$ws.Range('D17').Value = '''89.01'
$ws.Range('E17').Value = '  -0.08%  '
$ws.Range('D18').Value = '''0.06577'
$ws.Range('E18').Value = '  -2.25%  '
$ws.Range('D19').Value = '''0.00001018'
$ws.Range('E19').Value = '  -2.76%  '
$ws.Range('D20').Value = '''17.07'
$ws.Range('E20').Value = '  -1.22%  '
$ws.Range('E21').Value = '  -0.61%  '
$ws.Range('D22').Value = '27.340.58'
$ws.Range('E22').Value = '  -1.48%  '
$ws.Range('D23').Value = '''5.276'
$ws.Range('E23').Value = '  -3.73%  '
$ws.Range('D24').Value = '''10.77'
$ws.Range('E24').Value = '  -1.95%  '
$ws.Range('E25').Value = '  -1.69%  '
$ws.Range('D26').Value = '2.071.55'
$ws.Range('E26').Value = '  -1.73%  '
$ws.Range('D27').Value = '''155.63'
$ws.Range('E27').Value = '  -2.06%  '
$ws.Range('D28').Value = '''19.24'
$ws.Range('E28').Value = '  -2.69%  '
$ws.Range('D29').Value = '''2.034'
$ws.Range('E29').Value = '  -5.37%  '
$ws.Range('D30').Value = '''5.223'
$ws.Range('E30').Value = '  -4.66%  '
$ws.Range('D31').Value = '''117.44'
$ws.Range('E31').Value = '  -3.53%  '
$ws.Range('D32').Value = '''0.09274'
$ws.Range('E32').Value = '  -1.79%  '
$ws.Range('D33').Value = '''0.9278'
$ws.Range('E33').Value = '  -5.45%  '
$ws.Range('D34').Value = '''3.564'
$ws.Range('E34').Value = '  -1.79%  '
$ws.Range('D35').Value = '''5.202'
$ws.Range('E35').Value = '  -2.31%  '
$ws.Range('D36').Value = '''1.306'
$ws.Range('E36').Value = '  -3.07%  '
$ws.Range('D37').Value = '''0.05893'
$ws.Range('E37').Value = '  -2.82%  '
$ws.Range('D38').Value = '''0.02165'
$ws.Range('E38').Value = '  -3.27%  '
$ws.Range('B39').Value = 'Frax'
$ws.Range('C39').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D39').Value = '''1.002'
$ws.Range('E39').Value = '  -0.78%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = '''8.033'
$ws.Range('E40').Value = '  -3.74%  '
$ws.Range('D41').Value = '''1.134'
$ws.Range('E41').Value = '  -5.44%  '
$ws.Range('D42').Value = '''0.5716'
$ws.Range('E42').Value = '  -4.28%  '
$ws.Range('D43').Value = '''0.1811'
$ws.Range('E43').Value = '  -4.29%  '
$ws.Range('D44').Value = '''9.873'
$ws.Range('E44').Value = '  -4.77%  '
$ws.Range('D45').Value = '''1.279'
$ws.Range('E45').Value = '  +2.61%  '
$ws.Range('D46').Value = '''11.82'
$ws.Range('E46').Value = '  -3.71%  '
$ws.Range('D47').Value = '''0.5366'
$ws.Range('E47').Value = '  -4.95%  '
$ws.Range('D48').Value = '''1.860'
$ws.Range('E48').Value = '  -3.88%  '
$ws.Range('D49').Value = '''0.06563'
$ws.Range('E49').Value = '  -3.24%  '
$ws.Range('D50').Value = '''109.31'
$ws.Range('E50').Value = '  -2.28%  '
$ws.Range('D51').Value = '''1.003'
$ws.Range('E51').Value = '  -33.97%  '
